# "rebuilding the stats table in english"
# Replace the French/mixed day-title captions in column G (rows 2-44) of
# the GDMBR stats sheet with the new numbered English titles, adjust the
# sheet selection/viewport and widen column G to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: new English day titles (row 2 = day 43 ... row 44 = day 1) ---
$ws.Range("G2").Value  = "43. Sunset Post"
$ws.Range("G3").Value  = "42. Antelope Wells: the final stretch 💪. Really?"
$ws.Range("G4").Value  = "41. Hachita"
$ws.Range("G5").Value  = "39. Camaraderie"
$ws.Range("G6").Value  = "39. Beaverhead Ranch"
$ws.Range("G7").Value  = "38. Pie Town and Toaster House"
$ws.Range("G8").Value  = "37. Grants"
$ws.Range("G9").Value  = "36. Deserticus"
$ws.Range("G10").Value = "35. Abiquiu to Cuba"
$ws.Range("G11").Value = "34. Abiquiu Lodge"
$ws.Range("G12").Value = "33. Land of Enchantment - Really?"
$ws.Range("G13").Value = "32. Altus Maximus"
$ws.Range("G14").Value = "31. Del Norte - The Wild Colorado"
$ws.Range("G15").Value = "30. Cow-Girl encounter"
$ws.Range("G16").Value = "29. Salida and a Bit of Rest"
$ws.Range("G17").Value = "28. Boreas Pass 3,500 m"
$ws.Range("G18").Value = "27. Colorado Has Relief!"
$ws.Range("G19").Value = "26. Leaving Brush Mountain Lodge"
$ws.Range("G20").Value = "25. Hello Colorado!"
$ws.Range("G21").Value = "24. Things Go Wrong in Rawlins"
$ws.Range("G22").Value = "23. The Great Basin #2"
$ws.Range("G23").Value = "22. The Great Basin"
$ws.Range("G24").Value = "21. Pinendale"
$ws.Range("G25").Value = "20. Togwotee and Union Pass"
$ws.Range("G26").Value = "19. Welcome to Wyoming"
$ws.Range("G27").Value = "18. Old Oregon Short Line"
$ws.Range("G28").Value = "17. Today is a Good Day"
$ws.Range("G29").Value = "16. Storms"
$ws.Range("G30").Value = "15. Fleecer Ridge"
$ws.Range("G31").Value = "14. It's been 2 Weeks"
$ws.Range("G32").Value = "13. Lava Mountain and Butte"
$ws.Range("G33").Value = "12. Helena and Park Lake"
$ws.Range("G34").Value = "11. The 1,000 km mark"
$ws.Range("G35").Value = "10. The Alpaca Farm"
$ws.Range("G36").Value = "9. Overcoming 2075mD+"
$ws.Range("G37").Value = "8. Camping at Swan!"
$ws.Range("G38").Value = "7. A week already!"
$ws.Range("G39").Value = "6. The 49th Parallel"
$ws.Range("G40").Value = "5. The Titan"
$ws.Range("G41").Value = "4. The King"
$ws.Range("G42").Value = "3. The Day of the Grand Depart"
$ws.Range("G43").Value = "2. Arrival in Banff"
$ws.Range("G44").Value = "1. Not Showing Off on the First Day"

# --- widen column G so the longer English captions fit ---
$ws.Columns("G").ColumnWidth = 38.83

# --- update the viewport: scroll back to top-left, move selection to G59 ---
$ws.Activate()
$ws.Range("G59").Select() | Out-Null
